$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/fixed Cypher query for the "Participant ID" tab (B2) - corrects the
# failing instrument-model filter test case by re-scoping the participant
# match before collecting samples.
$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.instrument_model in ['Illumina HiSeq 2500']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id Limit 100
"@

$ws.Range("B2").Value = $newQuery

# Row 2 needs to be taller to fit the longer replacement query text.
$ws.Rows.Item(2).RowHeight = 300

# Move the active selection from C2 to B2 (also clears the stale
# top-left-cell scroll position that pointed at row 2).
$ws.Range("B2").Select() | Out-Null
